# Insert a new weekly record at row 275, shifting the existing rows 275:342 down to 276:343.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(275).Insert()

$ws.Range("A275").Value = 10
$ws.Range("B275").Value = "Vega Modelo de Temuco"
$ws.Range("C275").Value = "La Araucanía"
$ws.Range("D275").Value = 44754
$ws.Range("E275").Value = 9
$ws.Range("F275").Value = 100112009
$ws.Range("G275").Value = "Acelga"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 30
$ws.Range("K275").Value = 8000
$ws.Range("L275").Value = 8000
$ws.Range("M275").Value = 8000
$ws.Range("N275").Value = "$/docena de atados (12 kilos)"
$ws.Range("O275").Value = "Región Metropolitana"
$ws.Range("P275").Value = 667
$ws.Range("Q275").Value = 12
$ws.Range("R275").Value = "Hortaliza"
